$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Coin name, Link, Price (D), Volume(1h) (E), IsNumericPrice
$rows = @(
    @('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '27.513.91', '  +4.38%  ', $false),
    @('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.842.45', '  +3.76%  ', $false),
    @('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.032', '  +2.99%  ', $true),
    @('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '318.27', '  +4.18%  ', $true),
    @('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.027', '  +2.69%  ', $true),
    @('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4374', '  +3.38%  ', $true),
    @('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3732', '  +3.84%  ', $true),
    @('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07388', '  +3.35%  ', $true),
    @('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.8750', '  +4.71%  ', $true),
    @('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '21.42', '  +5.05%  ', $true),
    @('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.853.27', '  +4.09%  ', $false),
    @('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.480', '  +4.56%  ', $true),
    @('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.682', '  +3.80%  ', $true),
    @('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07157', '  +4.22%  ', $true),
    @('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '82.61', '  +4.55%  ', $true),
    @('BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.033', '  +2.75%  ', $true),
    @('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000009022', '  +4.41%  ', $true),
    @('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.027', '  +2.62%  ', $true),
    @('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '15.41', '  +3.55%  ', $true),
    @('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '27.534.04', '  +4.44%  ', $false),
    @('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.226', '  +3.04%  ', $true),
    @('Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '11.30', '  +3.84%  ', $true),
    @('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.065.79', '  +3.58%  ', $false),
    @('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '157.09', '  +3.62%  ', $true),
    @('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.913', '  +7.03%  ', $true),
    @('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '18.65', '  +3.63%  ', $true),
    @('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '5.267', '  +4.01%  ', $true),
    @('LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.929', '  +5.56%  ', $true),
    @('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '116.22', '  +1.62%  ', $true),
    @('Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.09069', '  +2.92%  ', $true),
    @('ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.207', '  +7.85%  ', $true),
    @('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7637', '  +5.42%  ', $true),
    @('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.486', '  +3.91%  ', $true),
    @('HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.877', '  +5.40%  ', $true),
    @('Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '1.030', '  +2.92%  ', $true),
    @('TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.147', '  +6.18%  ', $true),
    @('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01967', '  +4.70%  ', $true),
    @('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05252', '  +2.70%  ', $true),
    @('TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.5172', '  +5.27%  ', $true),
    @('MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.786', '  +7.00%  ', $true),
    @('Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1664', '  +3.79%  ', $true),
    @('FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.561', '  +4.20%  ', $true),
    @('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '8.499', '  +6.74%  ', $true),
    @('PaxosStandard', 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax', '1.030', '  +2.83%  ', $true),
    @('Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '109.14', '  +4.69%  ', $true),
    @('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '10.61', '  +4.52%  ', $true),
    @('PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.030', '  +3.00%  ', $true),
    @('NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.701', '  +4.33%  ', $true),
    @('Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.4637', '  +4.44%  ', $true),
    @('RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.902', '  +10.72%  ', $true)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]

    $priceCell = $ws.Cells.Item($r, 4)
    if ($row[4]) {
        # Price string looks numeric (e.g. "1.032") - force text storage
        # so it is preserved exactly like the source sheet (inline string),
        # instead of Excel auto-converting it to a number.
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $row[2]
        $priceCell.Style = "Normal"
    } else {
        $priceCell.Value = $row[2]
    }

    $ws.Cells.Item($r, 5).Value = $row[3]
    $r = $r + 1
}
